# Iteration v0.7 -- Split Camel Case hashtags
$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")
$wsSteps   = $wb.Worksheets.Item("Steps")

# --- Results sheet: bump the v0.2 "Test" column header to v0.7, and add
#     the new iteration's Accuracy / FPR / F1 scores in column H ---
$wsResults.Range("H1").Value = 0.7

$wsResults.Range("H10").Value = 0.91969696969697
$wsResults.Range("H11").Value = 0.0533070088845015
$wsResults.Range("H12").Value = 0.85778175313059

# --- Steps sheet: document the new step ---
$wsSteps.Range("A21").Value = 0.7
$wsSteps.Range("B21").Value = "split Camel case hashtags"

# --- restore/update each sheet's own selection cursor; select Results
#     first so that re-selecting on Steps leaves Steps as the active tab,
#     matching the workbook's original active-sheet state ---
$wsResults.Range("H12").Select()
$wsSteps.Range("B21").Select()
